$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.5704150199890137
$ws.Range("D3").Value = 0.1509251594543457
$ws.Range("D4").Value = 2.053256988525391
$ws.Range("D5").Value = 0.04685783386230469
$ws.Range("D6").Value = 9.754637956619263
$ws.Range("D7").Value = 0.01299595832824707
$ws.Range("D8").Value = 0.1078310012817383
$ws.Range("D9").Value = 1.276360988616943
$ws.Range("D10").Value = 5.714885950088501
$ws.Range("D11").Value = 0.3073081970214844
$ws.Range("D12").Value = 0.01186394691467285
$ws.Range("D13").Value = 0.03953409194946289
$ws.Range("D14").Value = 14.33312797546387
$ws.Range("D15").Value = 99.24435091018677
$ws.Range("D16").Value = 2.313521146774292
$ws.Range("D17").Value = 0.01349592208862305
$ws.Range("D18").Value = 0.08600306510925293
$ws.Range("D19").Value = 0.4214069843292236
$ws.Range("D20").Value = 0.6422049999237061
$ws.Range("D21").Value = 3.964146852493286
$ws.Range("D22").Value = 0.09005188941955566
$ws.Range("D23").Value = 0.01261591911315918
$ws.Range("D24").Value = 0.03083395957946777
$ws.Range("D25").Value = 0.02793502807617188
$ws.Range("D26").Value = 0.0130620002746582
$ws.Range("D27").Value = 0.06688213348388672
$ws.Range("D28").Value = 0.01235103607177734
$ws.Range("D29").Value = 0.1533589363098145
$ws.Range("D30").Value = 2.0609290599823
$ws.Range("D31").Value = 0.493441104888916
$ws.Range("D32").Value = 0.03567218780517578
$ws.Range("D33").Value = 9.741075038909912
$ws.Range("D34").Value = 0.01178598403930664
$ws.Range("D35").Value = 0.1100990772247314
$ws.Range("D36").Value = 0.04321789741516113
$ws.Range("D37").Value = 1.743595838546753
$ws.Range("D38").Value = 0.4172549247741699
$ws.Range("D39").Value = 25.06800580024719
$ws.Range("D40").Value = 179.1279811859131
$ws.Range("D41").Value = 1.26123309135437
$ws.Range("D42").Value = 0.1114749908447266
$ws.Range("D43").Value = 6.47984504699707
$ws.Range("D44").Value = 0.04830098152160645
$ws.Range("D45").Value = 0.406527042388916
$ws.Range("D46").Value = 0.842940092086792
$ws.Range("D47").Value = 2.003684997558594
$ws.Range("D48").Value = 0.7498798370361328
$ws.Range("D49").Value = 0.03198504447937012
$ws.Range("D50").Value = 0.1816260814666748
$ws.Range("D51").Value = 0.01298999786376953
$ws.Range("D52").Value = 0.05528092384338379
$ws.Range("D53").Value = 0.2879509925842285
$ws.Range("D54").Value = 0.01412296295166016
$ws.Range("D55").Value = 0.07988715171813965
$ws.Range("D56").Value = 0.0333409309387207
$ws.Range("D57").Value = 5.730488061904907
$ws.Range("D58").Value = 0.2275040149688721
$ws.Range("D59").Value = 53.30519104003906
$ws.Range("D60").Value = 0.0175929069519043
$ws.Range("D61").Value = 0.06056118011474609
$ws.Range("D62").Value = 3.676536083221436
$ws.Range("D63").Value = 0.5459098815917969
$ws.Range("D64").Value = 1.279531002044678
$ws.Range("D65").Value = 0.01209402084350586
$ws.Range("D66").Value = 0.04854607582092285
$ws.Range("D67").Value = 0.1074941158294678
$ws.Range("D68").Value = 0.3114171028137207
$ws.Range("D69").Value = 5.765976905822754
$ws.Range("D70").Value = 0.01941585540771484
$ws.Range("D71").Value = 0.3508529663085938
$ws.Range("D72").Value = 1.050707101821899
$ws.Range("D73").Value = 0.08148503303527832
$ws.Range("D74").Value = 0.1732320785522461
